$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column AY: header "pr1el" plus a numeric value per row (rows that
# already have #N/A data past column AN are skipped, matching the source).
$ws.Range("AY1").Value = "pr1el"

$ay = @{
    2  = 0.20133999999999999
    3  = 0.18182000000000001
    4  = 0.25502999999999998
    5  = 0.25925999999999999
    6  = 0.12739
    8  = 0.55769000000000002
    9  = 0.44374999999999998
    10 = 0.10256
    11 = 0.24285999999999999
    13 = 0.2
    14 = 0.12820999999999999
    15 = 0.39285999999999999
    16 = 0.35503000000000001
    17 = 0.32184000000000001
    18 = 0.31172
    21 = 0.16667000000000001
}

foreach ($row in $ay.Keys) {
    $ws.Cells.Item($row, 51).Value = $ay[$row]
}

# Restore the active selection to reflect the newly written corner cell.
$ws.Range("AY22").Select()
